$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = 'Cluster Name'
$ws.Cells.Item(1, 2).Value = 'Activecases'

# Data rows
$ws.Cells.Item(2, 1).Value = '3398 BlueCross Elly Kay Mordialloc'
$ws.Cells.Item(2, 2).Value = 29
$ws.Cells.Item(3, 1).Value = '3749 Rosebrook - McKenzie Aged Care Rosebud'
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(4, 1).Value = '4257 BlueCross The Gables Camberwell'
$ws.Cells.Item(4, 2).Value = 16
$ws.Cells.Item(5, 1).Value = '44404 Castlemaine North Primary School Castlemaine'
$ws.Cells.Item(5, 2).Value = 14
$ws.Cells.Item(6, 1).Value = '44622 Grey Street Primary School Traralgon'
$ws.Cells.Item(6, 2).Value = 12
$ws.Cells.Item(7, 1).Value = '44642 Irymple South Primary School Irymple South'
$ws.Cells.Item(7, 2).Value = 13
$ws.Cells.Item(8, 1).Value = '4479 Whittlesea Lodge Whittlesea'
$ws.Cells.Item(8, 2).Value = 15
$ws.Cells.Item(9, 1).Value = '45168 Ranfurly Primary School Mildura'
$ws.Cells.Item(9, 2).Value = 13
$ws.Cells.Item(10, 1).Value = '45275 Lalor Gardens Primary School Lalor'
$ws.Cells.Item(10, 2).Value = 11
$ws.Cells.Item(11, 1).Value = '52390 Our Lady of the Way Catholic Primary School Wallan'
$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(12, 1).Value = '52777 Mirripoa Primary School Mount Duneed School Camp'
$ws.Cells.Item(12, 2).Value = 21
$ws.Cells.Item(13, 1).Value = 'Alfred Health The Alfred Hospital Melbourne'
$ws.Cells.Item(13, 2).Value = 12
$ws.Cells.Item(14, 1).Value = 'Confirmed Omicron Sircuit Bar Fitzroy'
$ws.Cells.Item(14, 2).Value = 14
$ws.Cells.Item(15, 1).Value = 'Confirmed Omicron Variant The Peel Hotel Collingwood'
$ws.Cells.Item(15, 2).Value = 13
$ws.Cells.Item(16, 1).Value = 'Feathertop Chalet Harrietville'
$ws.Cells.Item(16, 2).Value = 18
$ws.Cells.Item(17, 1).Value = 'St Pauls Cathedral'
$ws.Cells.Item(17, 2).Value = 62
$ws.Cells.Item(18, 1).Value = 'St Vincents Hospital Melbourne Emergency Department Fitzroy'
$ws.Cells.Item(18, 2).Value = 19
$ws.Cells.Item(19, 1).Value = 'St. Vincent''s Hospital Melbourne Fitzroy'
$ws.Cells.Item(19, 2).Value = 19
$ws.Cells.Item(20, 1).Value = 'The Emerson Rooftop Bar and Club South Yarra'
$ws.Cells.Item(20, 2).Value = 14
$ws.Cells.Item(21, 1).Value = 'The Hatter and the Hare Bayswater'
$ws.Cells.Item(21, 2).Value = 16
$ws.Cells.Item(22, 1).Value = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville'
$ws.Cells.Item(22, 2).Value = 11
